# Edit script for "Rapport continue de stage.docx"
# 1) Merge the three runs (with gramStart/gramEnd proofErr markers) that make up
#    the "Accueil de l'equipe (...)" bullet into a single plain run - a
#    same-text Find/Replace collapses the run/proofErr boundaries for us.
# 2) Append the new "13/06" journal section (10 new bullet paragraphs) after the
#    last paragraph of the document, via Range.InsertXML so the exact run /
#    proofErr layout from the authored OOXML is reproduced.

$d = $word.ActiveDocument

# --- Change 1: tidy the "Accueil de l'equipe" bullet -----------------------
$accueilText = "Accueil de l’équipe (briefing, compte rendu des taches à faire…)"
$d.Content.Find.Execute($accueilText, $true, $false, $false, $false, $false, $true, 1, $false, $accueilText, 2) | Out-Null

# --- Change 2: append the new 13/06 section at the end of the document -----
$endPos = $d.Content.End
$insertionPoint = $d.Range($endPos, $endPos)
$newSectionXml = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="Paragraphedeliste"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve">Faits : </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Paragraphedeliste"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>Code réorganisé</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Paragraphedeliste"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve">Page contact fini avec Google </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>forms</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve"> (envoie de mail lors de formulaire soumis)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Paragraphedeliste"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve">Avancement du </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>footer</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve"> avec mini </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>map</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve"> de l’adresse et une rubrique « Navigation »</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Paragraphedeliste"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>13/06 :</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Paragraphedeliste"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve">But : </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Paragraphedeliste"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>Faire une autre page (que ce soit formulaire pour proposition de formation ou page de présentation)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Paragraphedeliste"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve">Finir le </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>footer</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve"> avec les badges et conditions avec modalités </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>ect</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>…</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Paragraphedeliste"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve">Page home : </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>acces</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve"> au formulaire de </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>prop</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve"> de formation </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>necessaire</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Paragraphedeliste"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Changement du </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t>forms</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="26"/>
          <w:szCs w:val="26"/>
        </w:rPr>
        <w:t xml:space="preserve"> concernant l’adresse mail réceptrice (création d’un nouveau compte google pour ça)</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($newSectionXml)
